$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows right before existing row 631, shifting the old
# rows 631-646 down to 635-650 (dimension grows from R646 to R650).
$ws.Rows("631:634").Insert()

# --- New row 631: Conconina(o) ---
$ws.Range("A631").Value2 = 10
$ws.Range("B631").Value2 = "Vega Modelo de Temuco"
$ws.Range("C631").Value2 = "La Araucanía"
$ws.Range("D631").Value2 = 44448
$ws.Range("E631").Value2 = 9
$ws.Range("F631").Value2 = 100112033
$ws.Range("G631").Value2 = "Lechuga"
$ws.Range("H631").Value2 = "Conconina(o)"
$ws.Range("I631").Value2 = "Primera"
$ws.Range("J631").Value2 = 185
$ws.Range("K631").Value2 = 7000
$ws.Range("L631").Value2 = 9000
$ws.Range("M631").Value2 = 8351
$ws.Range("N631").Value2 = "$/caja 10 unidades"
$ws.Range("O631").Value2 = "Región Metropolitana"
$ws.Range("P631").Value2 = 835
$ws.Range("Q631").Value2 = 10
$ws.Range("R631").Value2 = "Hortaliza"

# --- New row 632: Escarola ---
$ws.Range("A632").Value2 = 10
$ws.Range("B632").Value2 = "Vega Modelo de Temuco"
$ws.Range("C632").Value2 = "La Araucanía"
$ws.Range("D632").Value2 = 44448
$ws.Range("E632").Value2 = 9
$ws.Range("F632").Value2 = 100112033
$ws.Range("G632").Value2 = "Lechuga"
$ws.Range("H632").Value2 = "Escarola"
$ws.Range("I632").Value2 = "Primera"
$ws.Range("J632").Value2 = 1580
$ws.Range("K632").Value2 = 11000
$ws.Range("L632").Value2 = 13000
$ws.Range("M632").Value2 = 11620
$ws.Range("N632").Value2 = "$/caja 15 unidades"
$ws.Range("O632").Value2 = "Provincia del Elquí"
$ws.Range("P632").Value2 = 775
$ws.Range("Q632").Value2 = 15
$ws.Range("R632").Value2 = "Hortaliza"

# --- New row 633: Francesa morada ---
$ws.Range("A633").Value2 = 10
$ws.Range("B633").Value2 = "Vega Modelo de Temuco"
$ws.Range("C633").Value2 = "La Araucanía"
$ws.Range("D633").Value2 = 44448
$ws.Range("E633").Value2 = 9
$ws.Range("F633").Value2 = 100112033
$ws.Range("G633").Value2 = "Lechuga"
$ws.Range("H633").Value2 = "Francesa morada"
$ws.Range("I633").Value2 = "Primera"
$ws.Range("J633").Value2 = 125
$ws.Range("K633").Value2 = 7000
$ws.Range("L633").Value2 = 7000
$ws.Range("M633").Value2 = 7000
$ws.Range("N633").Value2 = "$/caja 15 unidades"
$ws.Range("O633").Value2 = "Región Metropolitana"
$ws.Range("P633").Value2 = 467
$ws.Range("Q633").Value2 = 15
$ws.Range("R633").Value2 = "Hortaliza"

# --- New row 634: Marina ---
$ws.Range("A634").Value2 = 10
$ws.Range("B634").Value2 = "Vega Modelo de Temuco"
$ws.Range("C634").Value2 = "La Araucanía"
$ws.Range("D634").Value2 = 44448
$ws.Range("E634").Value2 = 9
$ws.Range("F634").Value2 = 100112033
$ws.Range("G634").Value2 = "Lechuga"
$ws.Range("H634").Value2 = "Marina"
$ws.Range("I634").Value2 = "Primera"
$ws.Range("J634").Value2 = 75
$ws.Range("K634").Value2 = 7000
$ws.Range("L634").Value2 = 7000
$ws.Range("M634").Value2 = 7000
$ws.Range("N634").Value2 = "$/caja 15 unidades"
$ws.Range("O634").Value2 = "Región del Maule"
$ws.Range("P634").Value2 = 467
$ws.Range("Q634").Value2 = 15
$ws.Range("R634").Value2 = "Hortaliza"
